# BISAGRA MUNICION DISMAY price list refresh:
# - bump the price-list date in A1 by one month
# - update the three unit prices in column D (rows 33-35)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 45436

$ws.Range("D33").Value = 1687.737
$ws.Range("D34").Value = 1275.478
$ws.Range("D35").Value = 949.728
